$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 4 through 8 first
$ws.Rows("4:8").Delete()

# Update A2 and A3 with the new combined values
$ws.Range("A2").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"
$ws.Range("A3").Value = "('Knight', ['Token Creature — Knight', 'Vigilance', '2/2'])"
